# Update Sweden / Brokerage & Investment Banking capital-structure rows.
# Rows 2-4 get refreshed metrics; former row 5 (FX International AB) moves to
# row 6 with new figures, and Nordnet AB (publ) is inserted as the new row 5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Cells.Item(2, 1).Value = "Sweden"
$c = $ws.Cells.Item(2, 2)
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$ws.Cells.Item(2, 3).Value = "Brokerage & Investment Banking"
$ws.Cells.Item(2, 4).Value = 0.133
$ws.Cells.Item(2, 5).Value = 0.22
$ws.Cells.Item(2, 7).Value = 0.01199221154088401
$ws.Cells.Item(2, 8).Value = -0.001288478911055365
$ws.Cells.Item(2, 9).Value = -0.001624995051233233
$ws.Cells.Item(2, 10).Value = -0.001399811901259589
$ws.Cells.Item(2, 11).Value = 223.17
$ws.Cells.Item(2, 12).Value = 0.4016059198047847
$ws.Cells.Item(2, 13).Value = 3689.881
$ws.Cells.Item(2, 14).Value = 0.4384359096531492
$ws.Cells.Item(2, 15).Value = 16.53394721512748
$ws.Cells.Item(2, 16).Value = 3689.881
$ws.Cells.Item(2, 17).Value = 0.4384359096531492
$ws.Cells.Item(2, 18).Value = 16.53394721512748
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 364.83
$ws.Cells.Item(2, 22).Value = 0.04334952073488506
$ws.Cells.Item(2, 23).Value = 0.2973182262797522
$ws.Cells.Item(2, 24).Value = 0.0291751072644959
$ws.Cells.Item(2, 25).Value = 0.2681431190152563
$ws.Cells.Item(2, 26).Value = 0.1377954376496809
$ws.Cells.Item(2, 27).Value = 0
$ws.Cells.Item(2, 28).Value = 0.02904909240069335
$ws.Cells.Item(2, 29).Value = -0.03202441581140694
$ws.Cells.Item(2, 30).Value = 5762.5
$ws.Cells.Item(2, 31).Value = 0
$ws.Cells.Item(2, 32).Value = 5762.5
$ws.Cells.Item(2, 33).Value = 5397.67
$ws.Cells.Item(2, 34).Value = 0.4064249346369964
$ws.Cells.Item(2, 35).Value = 0.8859373870958319
$ws.Cells.Item(2, 36).Value = 0.3907481568995373
$ws.Cells.Item(2, 37).Value = 0.8791594864795312
$ws.Cells.Item(2, 38).Value = 0.003
$ws.Cells.Item(2, 39).Value = 0.003
$ws.Cells.Item(2, 40).Value = -6402.777777777777
$ws.Cells.Item(2, 41).Value = -301
$ws.Cells.Item(2, 42).Value = -5997.411111111111
$ws.Cells.Item(2, 43).Value = -301

# --- Row 3 ---
$ws.Cells.Item(3, 1).Value = "Sweden"
$ws.Cells.Item(3, 2).Value = "Avanza Bank Holding AB (publ) (OM:AZA)"
$ws.Cells.Item(3, 3).Value = "Brokerage & Investment Banking"
$ws.Cells.Item(3, 4).Value = 0.186
$ws.Cells.Item(3, 5).Value = 0.232
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 115.9
$ws.Cells.Item(3, 12).Value = 0.4632294164668266
$ws.Cells.Item(3, 13).Value = 39.5
$ws.Cells.Item(3, 14).Value = 0.008985236914537885
$ws.Cells.Item(3, 15).Value = 0.3408110440034512
$ws.Cells.Item(3, 16).Value = 39.5
$ws.Cells.Item(3, 17).Value = 0.008985236914537885
$ws.Cells.Item(3, 18).Value = 0.3408110440034512
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 172.1
$ws.Cells.Item(3, 22).Value = 0.03914833602511317
$ws.Cells.Item(3, 23).Value = 0.6194548369855692
$ws.Cells.Item(3, 24).Value = 0.0290978107880785
$ws.Cells.Item(3, 25).Value = 0.5903570261974908
$ws.Cells.Item(3, 26).Value = -5.900943396226414
$ws.Cells.Item(3, 27).Value = -0
$ws.Cells.Item(3, 28).Value = 0.02899753309667824
$ws.Cells.Item(3, 29).Value = -0.02899753309667824
$ws.Cells.Item(3, 30).Value = 31.7
$ws.Cells.Item(3, 31).Value = 0
$ws.Cells.Item(3, 32).Value = 31.7
$ws.Cells.Item(3, 33).Value = -140.4
$ws.Cells.Item(3, 34).Value = 0.007159311622024481
$ws.Cells.Item(3, 35).Value = 0.09789993823347745
$ws.Cells.Item(3, 36).Value = -0.03299104730126654
$ws.Cells.Item(3, 37).Value = -0.9255108767303889
$ws.Cells.Item(3, 38).Value = 0
$ws.Cells.Item(3, 39).Value = 0
$ws.Cells.Item(3, 40).ClearContents()
$ws.Cells.Item(3, 42).ClearContents()

# --- Row 4 ---
$ws.Cells.Item(4, 1).Value = "Sweden"
$ws.Cells.Item(4, 2).Value = "Mangold Fondkommission AB (OM:MANG)"
$ws.Cells.Item(4, 3).Value = "Brokerage & Investment Banking"
$ws.Cells.Item(4, 4).Value = 0.035
$ws.Cells.Item(4, 5).Value = 0.123
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 2.5
$ws.Cells.Item(4, 12).Value = 0.1445086705202312
$ws.Cells.Item(4, 13).Value = 0.381
$ws.Cells.Item(4, 14).Value = 0.004369266055045872
$ws.Cells.Item(4, 15).Value = 0.1524
$ws.Cells.Item(4, 16).Value = 0.381
$ws.Cells.Item(4, 17).Value = 0.004369266055045872
$ws.Cells.Item(4, 18).Value = 0.1524
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 0
$ws.Cells.Item(4, 22).Value = 0
$ws.Cells.Item(4, 23).Value = 0.25
$ws.Cells.Item(4, 24).Value = 0.02925240374091331
$ws.Cells.Item(4, 25).Value = 0.2207475962590867
$ws.Cells.Item(4, 26).Value = 1.251808972503618
$ws.Cells.Item(4, 27).Value = 0
$ws.Cells.Item(4, 28).Value = 0.02910065170470846
$ws.Cells.Item(4, 29).Value = -0.02910065170470846
$ws.Cells.Item(4, 30).Value = 1.5
$ws.Cells.Item(4, 31).Value = 0
$ws.Cells.Item(4, 32).Value = 1.5
$ws.Cells.Item(4, 33).Value = 1.5
$ws.Cells.Item(4, 34).Value = 0.01691093573844419
$ws.Cells.Item(4, 35).Value = 0.1013513513513513
$ws.Cells.Item(4, 36).Value = 0.01691093573844419
$ws.Cells.Item(4, 37).Value = 0.1013513513513513
$ws.Cells.Item(4, 38).Value = 0
$ws.Cells.Item(4, 39).Value = 0
$ws.Cells.Item(4, 40).ClearContents()
$ws.Cells.Item(4, 42).ClearContents()

# --- Row 5 ---
$ws.Cells.Item(5, 1).Value = "Sweden"
$ws.Cells.Item(5, 2).Value = "Nordnet AB (publ) (OM:SAVE)"
$ws.Cells.Item(5, 3).Value = "Brokerage & Investment Banking"
$ws.Cells.Item(5, 4).Value = 0.133
$ws.Cells.Item(5, 5).Value = 0.22
$ws.Cells.Item(5, 7).Value = 0.02558945908460472
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 105.7
$ws.Cells.Item(5, 12).Value = 0.3665048543689321
$ws.Cells.Item(5, 13).Value = 3650
$ws.Cells.Item(5, 14).Value = 0.9287768136593806
$ws.Cells.Item(5, 15).Value = 34.53169347209082
$ws.Cells.Item(5, 16).Value = 3650
$ws.Cells.Item(5, 17).Value = 0.9287768136593806
$ws.Cells.Item(5, 18).Value = 34.53169347209082
$ws.Cells.Item(5, 19).Value = 0
$ws.Cells.Item(5, 20).Value = 0
$ws.Cells.Item(5, 21).Value = 192.6
$ws.Cells.Item(5, 22).Value = 0.04900888063309498
$ws.Cells.Item(5, 23).Value = 0.3446364525595044
$ws.Cells.Item(5, 24).Value = 0.05154447442962004
$ws.Cells.Item(5, 25).Value = 0.2930919781298844
$ws.Cells.Item(5, 26).Value = 0.07101873968824644
$ws.Cells.Item(5, 27).Value = 0
$ws.Cells.Item(5, 28).Value = 0.03494817991810541
$ws.Cells.Item(5, 29).Value = -0.03494817991810541
$ws.Cells.Item(5, 30).Value = 5729.3
$ws.Cells.Item(5, 31).Value = 0
$ws.Cells.Item(5, 32).Value = 5729.3
$ws.Cells.Item(5, 33).Value = 5536.7
$ws.Cells.Item(5, 34).Value = 0.5931443597813483
$ws.Cells.Item(5, 35).Value = 0.9297038539553752
$ws.Cells.Item(5, 36).Value = 0.5848667948365833
$ws.Cells.Item(5, 37).Value = 0.9274359704517664
$ws.Cells.Item(5, 38).Value = 0
$ws.Cells.Item(5, 39).Value = 0
$ws.Cells.Item(5, 40).ClearContents()
$ws.Cells.Item(5, 42).ClearContents()
$ws.Cells.Item(5, 43).ClearContents()

# --- Row 6 ---
$ws.Cells.Item(6, 1).Value = "Sweden"
$ws.Cells.Item(6, 2).Value = "FX International AB (publ) (NGM:FXI)"
$ws.Cells.Item(6, 3).Value = "Brokerage & Investment Banking"
$ws.Cells.Item(6, 7).Value = 3.475728155339806
$ws.Cells.Item(6, 8).Value = 3.475728155339806
$ws.Cells.Item(6, 9).Value = 4.383495145631068
$ws.Cells.Item(6, 10).Value = 4.383495145631068
$ws.Cells.Item(6, 11).Value = -0.93
$ws.Cells.Item(6, 12).Value = 4.514563106796117
$ws.Cells.Item(6, 13).Value = -0
$ws.Cells.Item(6, 14).Value = -0
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = -0
$ws.Cells.Item(6, 17).Value = -0
$ws.Cells.Item(6, 18).Value = 0
$ws.Cells.Item(6, 19).Value = 0
$ws.Cells.Item(6, 21).Value = 0.13
$ws.Cells.Item(6, 22).Value = 0.04626334519572954
$ws.Cells.Item(6, 23).Value = -0.8017241379310346
$ws.Cells.Item(6, 24).Value = 0.02898623322517386
$ws.Cells.Item(6, 25).Value = -0.8307103711562085
$ws.Cells.Item(6, 26).Value = -0.4835680751173709
$ws.Cells.Item(6, 27).Value = -2.119718309859155
$ws.Cells.Item(6, 28).Value = 0.02898623322517386
$ws.Cells.Item(6, 29).Value = -2.148704543084329
$ws.Cells.Item(6, 30).Value = 0
$ws.Cells.Item(6, 31).Value = 0
$ws.Cells.Item(6, 32).Value = 0
$ws.Cells.Item(6, 33).Value = -0.13
$ws.Cells.Item(6, 34).Value = 0
$ws.Cells.Item(6, 35).Value = 0
$ws.Cells.Item(6, 36).Value = -0.04850746268656717
$ws.Cells.Item(6, 37).Value = -0.04088050314465409
$ws.Cells.Item(6, 38).Value = 0.003
$ws.Cells.Item(6, 39).Value = 0.003
$ws.Cells.Item(6, 40).Value = -0
$ws.Cells.Item(6, 41).Value = -301
$ws.Cells.Item(6, 42).Value = 0.1444444444444445
$ws.Cells.Item(6, 43).Value = -301
